# Updates parameter names (to match the companion params-mapping file) and
# replaces several "± X" / "± X%" placeholder text values in the
# Spike/Check Accuracy (and a few Field/Lab Duplicate) columns with actual
# numeric values. Cells that end up holding a genuine percentage (0.05,
# 0.15, ...) that previously held a "± N%" string get an explicit percentage
# number format so they still display as "5%"/"15%" etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Parameter name (column A) updates -------------------------------
$ws.Range("A2").Value  = "Water Temp"
$ws.Range("A4").Value  = "DO"
$ws.Range("A5").Value  = "DO"
$ws.Range("A6").Value  = "Conductivity"
$ws.Range("A7").Value  = "Conductivity"
$ws.Range("A12").Value = "Ortho P"
$ws.Range("A13").Value = "Ortho P"
$ws.Range("A17").Value = "Chl a"
$ws.Range("A18").Value = "Chl a"
$ws.Range("A19").Value = "E.coli"
$ws.Range("A20").Value = "E.coli"
$ws.Range("A21").Value = "E.coli"
$ws.Range("A22").Value = "E.coli"

# --- Replace "± ..." placeholder text with real numbers ---------------
$ws.Range("F3").Value  = 0.5
$ws.Range("G3").Value  = 0.5
$ws.Range("J3").Value  = 0.2

$ws.Range("J4").Value  = 0.05
$ws.Range("J5").Value  = 0.05

$ws.Range("J6").Value  = 50
$ws.Range("J7").Value  = 50

$ws.Range("F8").Value  = 1

$ws.Range("F10").Value = 0.02
$ws.Range("G10").Value = 0.01
$ws.Range("J10").Value = 0.01

$ws.Range("J11").Value = 0.15

$ws.Range("F12").Value = 0.01
$ws.Range("J12").Value = 0.01

$ws.Range("J13").Value = 0.15
$ws.Range("J14").Value = 0.15
$ws.Range("J15").Value = 0.15
$ws.Range("J16").Value = 0.15

$ws.Range("F17").Value = 2

# --- Percentage display for the cells that used to read "± N%" --------
$ws.Range("J4").NumberFormat  = "0%"
$ws.Range("J5").NumberFormat  = "0%"
$ws.Range("J11").NumberFormat = "0%"
$ws.Range("J13").NumberFormat = "0%"
$ws.Range("J14").NumberFormat = "0%"
$ws.Range("J15").NumberFormat = "0%"
$ws.Range("J16").NumberFormat = "0%"

# --- Misc view bookkeeping (best-effort; matches author's session) ----
$ws.Range("A22").Select()
